$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H11").Value = 18.11111
$ws.Range("I11").Value = 18.11111
$ws.Range("K11").Value = 18.11111
$ws.Range("M11").Value = 121.88889
$ws.Range("H32").Value = 2300.5625
$ws.Range("J32").Value = 2716.5
$ws.Range("L32").Value = 2716.5
$ws.Range("N32").Value = -3368.5
$ws.Range("H40").Value = 1886.375
$ws.Range("I40").Value = 1867.9231
$ws.Range("J40").Value = 1966.3334
$ws.Range("K40").Value = 1867.9231
$ws.Range("L40").Value = 1966.3334
$ws.Range("M40").Value = -1692.9231
$ws.Range("N40").Value = -2316.3334
$ws.Range("H42").Value = 3214.9167
$ws.Range("I42").Value = 1201.625
$ws.Range("J42").Value = 7241.5
$ws.Range("K42").Value = 3604.875
$ws.Range("L42").Value = 21724.5
$ws.Range("M42").Value = -3374.875
$ws.Range("N42").Value = -22184.5
$ws.Range("H43").Value = 0
$ws.Range("J43").Value = 0
$ws.Range("L43").Value = 0
$ws.Range("N43").ClearContents()
$ws.Range("H80").Value = 3800
$ws.Range("H83").Value = 3800
$ws.Range("H113").Value = 2048.3333
$ws.Range("I113").Value = 2040.2
$ws.Range("K113").Value = 2040.2
$ws.Range("M113").Value = 1213.8
$ws.Range("H114").Value = 69895
$ws.Range("J114").Value = 69895
$ws.Range("L114").Value = 69895
$ws.Range("N114").Value = -78573
$ws.Range("H132").Value = 1730.9445
$ws.Range("I132").Value = 1156.2941
$ws.Range("K132").Value = 3468.8823
$ws.Range("M132").Value = -938.8823000000002
$ws.Range("H137").Value = 1233
$ws.Range("I137").Value = 1243.2142
$ws.Range("J137").Value = 1090
$ws.Range("K137").Value = 3729.6426
$ws.Range("L137").Value = 3270
$ws.Range("M137").Value = -1179.6426
$ws.Range("N137").Value = -8370
$ws.Range("H138").Value = 3881.7837
$ws.Range("J138").Value = 5145.091
$ws.Range("L138").Value = 15435.273
$ws.Range("N138").Value = -25715.273
$ws.Range("H141").Value = 6498.25
$ws.Range("J141").Value = 6000
$ws.Range("L141").Value = 18000
$ws.Range("N141").Value = -28360

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2114.348
$ws.Range("I32").Value = 2067.6584
$ws.Range("K32").Value = 2067.6584
$ws.Range("M32").Value = -1780.6584
$ws.Range("H39").Value = 3638.3333
$ws.Range("I39").Value = 3638.3333
$ws.Range("K39").Value = 3638.3333
$ws.Range("M39").Value = -3118.3333
$ws.Range("H61").Value = 3100
$ws.Range("I61").Value = 3100
$ws.Range("K61").Value = 3100
$ws.Range("M61").Value = -2888
$ws.Range("H63").Value = 5106
$ws.Range("I63").Value = 4907.125
$ws.Range("K63").Value = 4907.125
$ws.Range("M63").Value = -4221.125
$ws.Range("H66").Value = 5106
$ws.Range("I66").Value = 4907.125
$ws.Range("K66").Value = 24535.625
$ws.Range("M66").Value = -21103.625
$ws.Range("H74").Value = 1122.4286
$ws.Range("I74").Value = 802
$ws.Range("J74").Value = 1699.2
$ws.Range("K74").Value = 802
$ws.Range("L74").Value = 1699.2
$ws.Range("M74").Value = 72
$ws.Range("N74").Value = -3447.2
$ws.Range("H77").Value = 1122.4286
$ws.Range("I77").Value = 802
$ws.Range("J77").Value = 1699.2
$ws.Range("K77").Value = 4010
$ws.Range("L77").Value = 8496
$ws.Range("M77").Value = 358
$ws.Range("N77").Value = -17232
$ws.Range("H132").Value = 2271.158
$ws.Range("I132").Value = 2266.9412
$ws.Range("J132").Value = 2307
$ws.Range("K132").Value = 6800.823600000001
$ws.Range("L132").Value = 6921
$ws.Range("M132").Value = -4270.823600000001
$ws.Range("N132").Value = -11981
$ws.Range("H136").Value = 3100
$ws.Range("I136").Value = 3100
$ws.Range("K136").Value = 9300
$ws.Range("M136").Value = -6750

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 250
$ws.Range("I11").Value = 4
$ws.Range("J11").Value = 332
$ws.Range("K11").Value = 4
$ws.Range("L11").Value = 332
$ws.Range("M11").Value = 136
$ws.Range("N11").Value = -612
$ws.Range("H64").Value = 963.75
$ws.Range("I64").Value = 965
$ws.Range("K64").Value = 965
$ws.Range("M64").Value = -740
$ws.Range("H67").Value = 963.75
$ws.Range("I67").Value = 965
$ws.Range("K67").Value = 965
$ws.Range("M67").Value = -185
$ws.Range("H94").Value = 639.0952
$ws.Range("I94").Value = 563.8125
$ws.Range("J94").Value = 880
$ws.Range("K94").Value = 563.8125
$ws.Range("L94").Value = 880
$ws.Range("M94").Value = -112.8125
$ws.Range("N94").Value = -1782
$ws.Range("H105").Value = 2981.6667
$ws.Range("I105").Value = 2998
$ws.Range("K105").Value = 2998
$ws.Range("M105").Value = -1251
$ws.Range("H107").Value = 1674.0938
$ws.Range("I107").Value = 1272.2693
$ws.Range("J107").Value = 3415.3333
$ws.Range("K107").Value = 1272.2693
$ws.Range("L107").Value = 3415.3333
$ws.Range("M107").Value = 647.7307000000001
$ws.Range("N107").Value = -7255.3333
$ws.Range("H130").Value = 0
$ws.Range("J130").Value = 0
$ws.Range("L130").Value = 0
$ws.Range("N130").ClearContents()

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H59").Value = 0
$ws.Range("J59").Value = 0
$ws.Range("L59").Value = 0
$ws.Range("N59").ClearContents()
$ws.Range("H107").Value = 1083.4
$ws.Range("I107").Value = 985.9286
$ws.Range("K107").Value = 985.9286
$ws.Range("M107").Value = 934.0714
$ws.Range("H135").Value = 28000
$ws.Range("J135").Value = 28000
$ws.Range("L135").Value = 28000
$ws.Range("N135").Value = -38140

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 8000
$ws.Range("I39").Value = 0
$ws.Range("J39").Value = 8000
$ws.Range("K39").Value = 0
$ws.Range("L39").Value = 24000
$ws.Range("M39").ClearContents()
$ws.Range("N39").Value = -24588
$ws.Range("H117").Value = 33503.668
$ws.Range("J117").Value = 50105.5
$ws.Range("L117").Value = 150316.5
$ws.Range("N117").Value = -157200.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 4798.6
$ws.Range("I132").Value = 4664.6665
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 13993.9995
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -11463.9995
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1637.25
$ws.Range("I22").Value = 1516.3334
$ws.Range("J22").Value = 2000
$ws.Range("K22").Value = 1516.3334
$ws.Range("L22").Value = 2000
$ws.Range("M22").Value = -1221.3334
$ws.Range("N22").Value = -2590
$ws.Range("H27").Value = 1637.25
$ws.Range("I27").Value = 1516.3334
$ws.Range("J27").Value = 2000
$ws.Range("K27").Value = 1516.3334
$ws.Range("L27").Value = 2000
$ws.Range("M27").Value = -1409.3334
$ws.Range("N27").Value = -2214
$ws.Range("H82").Value = 1042.6364
$ws.Range("I82").Value = 797.1429000000001
$ws.Range("J82").Value = 1472.25
$ws.Range("K82").Value = 797.1429000000001
$ws.Range("L82").Value = 1472.25
$ws.Range("M82").Value = -436.1429000000001
$ws.Range("N82").Value = -2194.25
$ws.Range("H85").Value = 1042.6364
$ws.Range("I85").Value = 797.1429000000001
$ws.Range("J85").Value = 1472.25
$ws.Range("K85").Value = 797.1429000000001
$ws.Range("L85").Value = 1472.25
$ws.Range("M85").Value = 450.8570999999999
$ws.Range("N85").Value = -3968.25
$ws.Range("H100").Value = 3749.182
$ws.Range("I100").Value = 3099.625
$ws.Range("K100").Value = 3099.625
$ws.Range("M100").Value = -2558.625
$ws.Range("H132").Value = 10456.167
$ws.Range("I132").Value = 13684.25
$ws.Range("K132").Value = 41052.75
$ws.Range("M132").Value = -38522.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 33749.75
$ws.Range("I14").Value = 34999.5
$ws.Range("J14").Value = 32500
$ws.Range("K14").Value = 34999.5
$ws.Range("L14").Value = 32500
$ws.Range("M14").Value = -34831.5
$ws.Range("N14").Value = -32836
$ws.Range("H132").Value = 8524.111000000001
$ws.Range("I132").Value = 12118.833
$ws.Range("K132").Value = 36356.499
$ws.Range("M132").Value = -33826.499
$ws.Range("H136").Value = 4308.6
$ws.Range("I136").Value = 4797.35
$ws.Range("J136").Value = 2353.6
$ws.Range("K136").Value = 14392.05
$ws.Range("L136").Value = 7060.799999999999
$ws.Range("M136").Value = -11842.05
$ws.Range("N136").Value = -12160.8
